# Insert a new weekly record for "Vega Modelo de Temuco - Durazno" above the
# existing row 314, shifting the remaining rows (314-331) down to (315-332).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(314).Insert()

$ws.Range("A314").Value = 10
$ws.Range("B314").Value = "Vega Modelo de Temuco"
$ws.Range("C314").Value = "La Araucanía"
$ws.Range("D314").Value = 44931
$ws.Range("E314").Value = 9
$ws.Range("F314").Value = "Fruta"
$ws.Range("G314").Value = 100103
$ws.Range("H314").Value = "Frutos de hueso (carozo)"
$ws.Range("I314").Value = 100103004
$ws.Range("J314").Value = "Durazno"
$ws.Range("K314").Value = "Doctor Davis"
$ws.Range("L314").Value = "Primera"
$ws.Range("M314").Value = 110
$ws.Range("N314").Value = 24000
$ws.Range("O314").Value = 24000
$ws.Range("P314").Value = 24000
$ws.Range("Q314").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R314").Value = "Región de O'Higgins"
$ws.Range("S314").Value = 1333
$ws.Range("T314").Value = 18
